$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197, shifting rows 197-278 down to 198-279
$ws.Rows.Item(197).Insert()

# Copy formatting for column D (date style) from the row above into the new row 197
$ws.Range("D196").Copy()
$ws.Range("D197").PasteSpecial(-4122) | Out-Null

# Populate the new row 197 with its values
$ws.Range("A197").Value = 4
$ws.Range("B197").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C197").Value = "Los Lagos"
$ws.Range("D197").Value = 44726
$ws.Range("E197").Value = 10
$ws.Range("F197").Value = 100112037
$ws.Range("G197").Value = "Cebollín"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 140
$ws.Range("K197").Value = 10000
$ws.Range("L197").Value = 10000
$ws.Range("M197").Value = 10000
$ws.Range("N197").Value = "$/paquete 36 unidades"
$ws.Range("O197").Value = "Región Metropolitana"
$ws.Range("P197").Value = 278
$ws.Range("Q197").Value = 36
$ws.Range("R197").Value = "Hortaliza"
